# Scheduled-runner data refresh: update Leve profit calculations
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 767.35297
$ws.Range("I38").Value = 75.07143000000001
$ws.Range("J38").Value = 3998
$ws.Range("K38").Value = 225.21429
$ws.Range("L38").Value = 11994
$ws.Range("M38").Value = 146.78571
$ws.Range("N38").Value = -12738
# Row 43
$ws.Range("H43").Value = 6924.067
$ws.Range("I43").Value = 6731.4287
$ws.Range("J43").Value = 7092.625
$ws.Range("K43").Value = 6731.4287
$ws.Range("L43").Value = 7092.625
$ws.Range("M43").Value = -6662.4287
$ws.Range("N43").Value = -7230.625
# Row 123
$ws.Range("H123").Value = 80000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 80000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 80000
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -89800
# Row 137
$ws.Range("H137").Value = 4285.4165
$ws.Range("I137").Value = 4545.2915
$ws.Range("J137").Value = 3765.6667
$ws.Range("K137").Value = 13635.8745
$ws.Range("L137").Value = 11297.0001
$ws.Range("M137").Value = -11085.8745
$ws.Range("N137").Value = -16397.0001
# Row 138
$ws.Range("H138").Value = 2393.8333
$ws.Range("I138").Value = 926.7368
$ws.Range("J138").Value = 2919.7737
$ws.Range("K138").Value = 2780.2104
$ws.Range("L138").Value = 8759.321100000001
$ws.Range("M138").Value = 2359.7896
$ws.Range("N138").Value = -19039.3211

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 537.5
$ws.Range("I4").Value = 556.6667
$ws.Range("J4").Value = 480
$ws.Range("K4").Value = 556.6667
$ws.Range("L4").Value = 480
$ws.Range("M4").Value = -440.6667
$ws.Range("N4").Value = -712
# Row 45
$ws.Range("H45").Value = 14708056
$ws.Range("I45").Value = 19232428
$ws.Range("J45").Value = 3847.375
$ws.Range("K45").Value = 19232428
$ws.Range("L45").Value = 3847.375
$ws.Range("M45").Value = -19232051
$ws.Range("N45").Value = -4601.375

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 1170.5714
$ws.Range("I22").Value = 896.6667
$ws.Range("K22").Value = 896.6667
$ws.Range("M22").Value = -723.6667
# Row 134
$ws.Range("H134").Value = 6670190
$ws.Range("I134").Value = 3527.077
$ws.Range("K134").Value = 10581.231
$ws.Range("M134").Value = -8046.231

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1016685.25
$ws.Range("I31").Value = 3459.6
$ws.Range("K31").Value = 3459.6
$ws.Range("M31").Value = -3164.6
# Row 34
$ws.Range("H34").Value = 1016685.25
$ws.Range("I34").Value = 3459.6
$ws.Range("K34").Value = 3459.6
$ws.Range("M34").Value = -3257.6
# Row 74
$ws.Range("H74").Value = 49924.75
$ws.Range("J74").Value = 49924.75
$ws.Range("L74").Value = 49924.75
$ws.Range("N74").Value = -51672.75
# Row 77
$ws.Range("H77").Value = 49924.75
$ws.Range("J77").Value = 49924.75
$ws.Range("L77").Value = 149774.25
$ws.Range("N77").Value = -158510.25
# Row 99
$ws.Range("H99").Value = 3669.3125
$ws.Range("I99").Value = 3522.1667
$ws.Range("J99").Value = 4110.75
$ws.Range("K99").Value = 3522.1667
$ws.Range("L99").Value = 4110.75
$ws.Range("M99").Value = -2024.1667
$ws.Range("N99").Value = -7106.75
# Row 126
$ws.Range("H126").Value = 3669.3125
$ws.Range("I126").Value = 3522.1667
$ws.Range("J126").Value = 4110.75
$ws.Range("K126").Value = 10566.5001
$ws.Range("L126").Value = 12332.25
$ws.Range("M126").Value = -8096.500100000001
$ws.Range("N126").Value = -17272.25
# Row 132
$ws.Range("H132").Value = 2579.9473
$ws.Range("I132").Value = 2295.2354
$ws.Range("K132").Value = 6885.706200000001
$ws.Range("M132").Value = -4355.706200000001

$ws = $wb.Worksheets.Item("CUL")
# Row 86
$ws.Range("H86").Value = 595.8
$ws.Range("I86").Value = 560
$ws.Range("J86").Value = 649.5
$ws.Range("K86").Value = 1680
$ws.Range("L86").Value = 1948.5
$ws.Range("M86").Value = -494
$ws.Range("N86").Value = -4320.5
# Row 89
$ws.Range("H89").Value = 595.8
$ws.Range("I89").Value = 560
$ws.Range("J89").Value = 649.5
$ws.Range("K89").Value = 5040
$ws.Range("L89").Value = 5845.5
$ws.Range("M89").Value = 888
$ws.Range("N89").Value = -17701.5
# Row 131
$ws.Range("H131").Value = 4626.3115
$ws.Range("J131").Value = 4626.3115
$ws.Range("L131").Value = 13878.9345
$ws.Range("N131").Value = -23958.9345

$ws = $wb.Worksheets.Item("GSM")
# Row 22
$ws.Range("H22").Value = 2338.6667
$ws.Range("J22").Value = 2000
$ws.Range("L22").Value = 2000
$ws.Range("N22").Value = -3058
# Row 97
$ws.Range("H97").Value = 1499.15
$ws.Range("I97").Value = 1533.2
$ws.Range("J97").Value = 1397
$ws.Range("K97").Value = 1533.2
$ws.Range("L97").Value = 1397
$ws.Range("M97").Value = -1037.2
$ws.Range("N97").Value = -2389
# Row 102
$ws.Range("H102").Value = 3380.55
$ws.Range("I102").Value = 2743.182
$ws.Range("K102").Value = 2743.182
$ws.Range("M102").Value = -1121.182
# Row 132
$ws.Range("H132").Value = 35722010
$ws.Range("I132").Value = 50003470
$ws.Range("K132").Value = 150010410
$ws.Range("M132").Value = -150007880

$ws = $wb.Worksheets.Item("LTW")
# Row 6
$ws.Range("H6").Value = 97000
$ws.Range("J6").Value = 97000
$ws.Range("L6").Value = 97000
$ws.Range("N6").Value = -97224
# Row 40
$ws.Range("H40").Value = 4908.3716
$ws.Range("I40").Value = 4164.273
$ws.Range("K40").Value = 4164.273
$ws.Range("M40").Value = -4028.273
# Row 46
$ws.Range("H46").Value = 3270.7368
$ws.Range("J46").Value = 4886.5713
$ws.Range("L46").Value = 4886.5713
$ws.Range("N46").Value = -5262.5713
# Row 96
$ws.Range("H96").Value = 99999
$ws.Range("J96").Value = 99999
$ws.Range("L96").Value = 99999
$ws.Range("N96").Value = -105491
# Row 132
$ws.Range("H132").Value = 1014221.8
$ws.Range("I132").Value = 23320.6
$ws.Range("J132").Value = 2005123
$ws.Range("K132").Value = 69961.79999999999
$ws.Range("L132").Value = 6015369
$ws.Range("M132").Value = -67431.79999999999
$ws.Range("N132").Value = -6020429

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 25007124
$ws.Range("J62").Value = 25007124
$ws.Range("L62").Value = 25007124
$ws.Range("N62").Value = -25008372
# Row 65
$ws.Range("H65").Value = 25007124
$ws.Range("J65").Value = 25007124
$ws.Range("L65").Value = 125035620
$ws.Range("N65").Value = -125041860
# Row 69
$ws.Range("H69").Value = 28090.334
$ws.Range("J69").Value = 28090.334
$ws.Range("L69").Value = 28090.334
$ws.Range("N69").Value = -29588.334
# Row 72
$ws.Range("H72").Value = 28090.334
$ws.Range("J72").Value = 28090.334
$ws.Range("L72").Value = 84271.00199999999
$ws.Range("N72").Value = -91759.00199999999
# Row 107
$ws.Range("H107").Value = 737.7931
$ws.Range("I107").Value = 732.7273
$ws.Range("K107").Value = 2198.1819
$ws.Range("M107").Value = -278.1819
